# Add new match records (fecha 2025-08-30, serial 45899) to the "Partidos" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")

$fecha = Get-Date -Year 2025 -Month 8 -Day 30 -Hour 0 -Minute 0 -Second 0

$rows = @(
    @("Jorge Gonzalez",             "Amarillo", "Arquero",       0, 0, $true,  3, 0, 0, 0, 0),
    @("Fabian Caicedo",             "Azul",     "Arquero",       0, 0, $true,  6, 1, 0, 0, 0),
    @("Edwing Yesid Castillo",      "Amarillo", "Mediocampista", 1, 0, $false, 0, 0, 0, 2, 0),
    @("Carlos Fernando Valencia",   "Amarillo", "Delantero",     3, 0, $false, 0, 0, 0, 0, 0),
    @("Armando Vieras",             "Amarillo", "Defensa",       1, 0, $false, 0, 0, 0, 0, 0),
    @("Juan Carlos Otero",          "Amarillo", "Mediocampista", 1, 0, $false, 0, 0, 0, 0, 0),
    @("Quintero ",                  "Amarillo", "Defensa",       0, 0, $false, 0, 0, 0, 1, 0),
    @("Juan David Espinal",         "Amarillo", "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @("Alexander Uribe",            "Azul",     "Mediocampista", 2, 0, $false, 0, 0, 0, 0, 0),
    @("Cesar Augusto Estrada",      "Azul",     "Delantero",     1, 0, $false, 0, 0, 0, 0, 0),
    @("Luis Carlos Arciniegas",     "Azul",     "Defensa",       0, 0, $false, 0, 0, 0, 1, 0),
    @("Invitado",                   "Amarillo", "Defensa",       0, 0, $false, 0, 0, 1, 0, 0),
    @("Luis David",                 "Amarillo", "Mediocampista", 0, 0, $false, 0, 0, 1, 0, 0)
)

$startRow = 403
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $fecha
    $ws.Cells.Item($r, 2).Value = $data[0]
    $ws.Cells.Item($r, 3).Value = $data[1]
    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]
    $ws.Cells.Item($r, 6).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[5]
    $ws.Cells.Item($r, 8).Value = $data[6]
    $ws.Cells.Item($r, 9).Value = $data[7]
    $ws.Cells.Item($r, 10).Value = $data[8]
    $ws.Cells.Item($r, 11).Value = $data[9]
    $ws.Cells.Item($r, 12).Value = $data[10]
}

# Scroll the frozen pane so row 401 is the first visible row below the
# header, then leave the selection on F418 (mirrors the author's on-save
# cursor position).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 401
$ws.Range("F418").Select()
